# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages have been generated:
#   - Status columns flip from "Ready for handoff" to
#     "Handed back: in sync with en-US" on every sheet that shows it.
#   - The per-language sheets gain their "Latest Target File" / "Latest
#     Handback File" filenames (with hyperlinks back to the source .md)
#     and a real "Latest Handback DateTime" instead of the 0001-01-01
#     placeholder.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1818a0fd51991cf95c9ea55f39c729d770529a5d/e2e/"

# File 1: 7f24a935-... , File 2: d8813d03-...
$file1Md  = "7f24a935-d87a-4a1a-a0db-bbb9c49e6697.md"
$file2Md  = "d8813d03-34b0-4a62-87af-6c99c4a0c316.md"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-language sheets: zh-cn and de-de
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; HandoffXlf1 = "7f24a935-d87a-4a1a-a0db-bbb9c49e6697.f4ccaf87d52c82aa73b706dc894cc2a5f2b76201.zh-cn.xlf"; HandoffXlf2 = "d8813d03-34b0-4a62-87af-6c99c4a0c316.d08284d261ab7cfcdaf1a5f89345109f97f94526.zh-cn.xlf"; HandbackTime = "2016-09-06 15:40:57" },
    @{ Sheet = "de-de"; HandoffXlf1 = "7f24a935-d87a-4a1a-a0db-bbb9c49e6697.f4ccaf87d52c82aa73b706dc894cc2a5f2b76201.de-de.xlf"; HandoffXlf2 = "d8813d03-34b0-4a62-87af-6c99c4a0c316.d08284d261ab7cfcdaf1a5f89345109f97f94526.de-de.xlf"; HandbackTime = "2016-09-06 15:41:22" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column (C) on both data rows
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Row 2 -> file 7f24a935-...
    $ws.Range("I2").Value = $file1Md
    $ws.Hyperlinks.Add($ws.Range("I2"), ($mdUrlBase + $file1Md), [Type]::Missing, [Type]::Missing, $file1Md) | Out-Null
    $ws.Range("J2").Value = $lang.HandoffXlf1
    $ws.Range("K2").Value = $lang.HandbackTime

    # Row 3 -> file d8813d03-...
    $ws.Range("I3").Value = $file2Md
    $ws.Hyperlinks.Add($ws.Range("I3"), ($mdUrlBase + $file2Md), [Type]::Missing, [Type]::Missing, $file2Md) | Out-Null
    $ws.Range("J3").Value = $lang.HandoffXlf2
    $ws.Range("K3").Value = $lang.HandbackTime

    # Widen columns so the new long values / status text are fully visible
    $ws.Range("C:C").ColumnWidth = 29.9777047293527
    $ws.Range("I:I").ColumnWidth = 40
    $ws.Range("J:J").ColumnWidth = 40
}

# Overview sheet: widen the zh-cn / de-de status columns (E, F) to fit the
# new, longer status text.
$overview.Range("E:E").ColumnWidth = 29.9777047293527
$overview.Range("F:F").ColumnWidth = 29.9777047293527
